$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the cells we are about to touch to stay text (avoids Excel
# auto-converting the numeric-looking strings into real numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

# Update Price column (D) for the rows whose value changed
$ws.Range("D2").Value = "242.57"
$ws.Range("D4").Value = "5.416"
$ws.Range("D5").Value = "0.05900"
$ws.Range("D6").Value = "3.444"
$ws.Range("D7").Value = "6.517"
$ws.Range("D8").Value = "0.8106"
$ws.Range("D9").Value = "0.9580"
$ws.Range("D10").Value = "0.01132"
$ws.Range("D11").Value = "0.1422"
$ws.Range("D12").Value = "0.07416"
$ws.Range("D13").Value = "0.03273"
$ws.Range("D14").Value = "0.03046"
$ws.Range("D15").Value = "0.09330"
$ws.Range("D16").Value = "3.850"
$ws.Range("D17").Value = "0.001576"
$ws.Range("D18").Value = "0.04670"
$ws.Range("D19").Value = "0.005875"
$ws.Range("D20").Value = "0.001265"
$ws.Range("D21").Value = "0.004896"
$ws.Range("D22").Value = "0.00006806"
$ws.Range("D23").Value = "3.585"
$ws.Range("D24").Value = "2.132"
$ws.Range("D26").Value = "0.1334"
$ws.Range("D27").Value = "0.0002286"
$ws.Range("D40").Value = "0.03954"
$ws.Range("D41").Value = "0.006188"
$ws.Range("D43").Value = "0.003003"
$ws.Range("D44").Value = "0.009934"
$ws.Range("D45").Value = "0.00005209"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").Value = "0.6666"
$ws.Range("D48").Value = "0.002383"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D50").Value = "0.0002002"

# Update Hora column (G) from 16 to 17 for all data rows (2-51)
for ($r = 2; $r -le 51; $r++) {
    $ws.Range("G$r").Value = "17"
}

